# Automatic update of files.
# - Truncate the fractional Ost/Nord coordinates on row 2 (Q2/R2) to integers.
# - Append a new observation record as row 3, mirroring the layout of row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 coordinates (Q2, R2) ---
$ws.Range("Q2").Value = 703230
$ws.Range("R2").Value = 7299262

# --- Append new row 3 ---
$ws.Range("A3").Value = 111867419
$ws.Range("B3").Value = 90658
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 4361
$ws.Range("F3").Value = "Orange taggsvamp"
$ws.Range("G3").Value = "Hydnellum aurantiacum"
$ws.Range("H3").Value = "(Batsch:Fr.) P.Karst."

$ws.Range("P3").Value = "Vikvallen, Pi lm"
$ws.Range("Q3").Value = 703160
$ws.Range("R3").Value = 7299375
$ws.Range("S3").Value = 10
$ws.Range("T3").Value = "Norrbotten"
$ws.Range("U3").Value = "Arvidsjaur"
$ws.Range("V3").Value = "Pite lappmark"
$ws.Range("W3").Value = "Arvidsjaur"

# Y3 / AA3 hold literal date-looking text ("2023-09-01"); force text storage
# (as row 2 already has) so Excel doesn't reinterpret them as date serials.
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2023-09-01"
$ws.Range("Z3").Value = "12:30"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2023-09-01"
$ws.Range("AB3").Value = "15:30"

$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

$ws.Range("AW3").Value = "Heike Kontermann"
$ws.Range("AX3").Value = "Heike Kontermann, Steve Daurer, Kirsten Stelling"
